$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vars = @(
  '_tejgct_r09gstcp',
  'tejgfun_f5ct06viv',
  '_tejgtotfun_f2opsegpc',
  '_tejgfun_f5ct05prots',
  '_devppimtotfun_f1trans',
  'devppimfun_f5r07ct05salud',
  '_tejgge_r09ct06acanf',
  'tejgfun_f5r07ct05cydep',
  'tdvgfun_f1ct06trans',
  'dfgpimpiafun_f1ct05prots',
  'tejgkft_rdetpc',
  '_tejgfun_f5r18ct06amb',
  '_tejgfun_f5r18ct05prots',
  'dfgdevpiagfun_f5r18ct06agropc',
  'devppimtotfun_f5r07salud',
  'pimgfun_f5r18ct05trans',
  'mod_110',
  '_tejgtotfun_f5r08opsegpc',
  'tdvgfun_f1ct06transpc',
  'pimgfun_f1ct06trans',
  'pimgct_r00gstcppc',
  '_tejgfun_f5ct06amb',
  '_pimgfun_f5r18ct06opseg',
  'devppimge_r00ct05biser',
  'tdvgfun_f5r18ct05cydep',
  '_devppimtotfun_f5r07trans',
  '_devppimfun_f5r18ct05salud',
  'empinc_2',
  'piagfun_f5r18ct05pgrcopc',
  'tejgfun_f5ct05trans',
  '_pimgfun_f5r07ct05salud',
  'tejgfun_f5r18ct05cydeppc',
  'devppimfun_f5r07ct05agro',
  '_pimgfun_f5r18ct06opsegpc',
  'mod_4',
  'dfgpimpiafun_f1ct06transpc',
  'dfgpimpiage_r08ct05dotra',
  'dfgpimpiafun_f5r08ct06agropc',
  'pimgfun_f2ct06san',
  'pimgfun_f1ct06san',
  '_devppimfun_f5r07ct05cydep',
  'compu_muni_1',
  'piagfun_f5r07ct05trans',
  'dfgdevpiagfun_f5r07ct05pgrcopc',
  '_pimgfun_f5ct06opsegpc',
  'mod_103',
  'devppimfun_f5r18ct06energ',
  '_devppimfun_f5ct05salud',
  '_dfgpimpiafun_f1ct06opsegpc',
  '_dfgdevpiagfun_f2ct05prsopc'
)

$scores = @(
  0.0254090279340744,
  0.0232990924268961,
  0.0223840344697237,
  0.01770706288516521,
  0.01620905660092831,
  0.01484724041074514,
  0.01294951979070902,
  0.01273930724710226,
  0.01232943683862686,
  0.01221377309411764,
  0.01138773560523987,
  0.01078783255070448,
  0.01048282720148563,
  0.009863639250397682,
  0.009600085206329823,
  0.009527770802378654,
  0.009338928386569023,
  0.008815799839794636,
  0.008676604367792606,
  0.008283902890980244,
  0.008233900181949139,
  0.007686045486479998,
  0.007376410998404026,
  0.007189224474132061,
  0.00667710741981864,
  0.006483918521553278,
  0.006331060081720352,
  0.006263700313866138,
  0.006196929141879082,
  0.006103032734245062,
  0.006074834149330854,
  0.005743281450122595,
  0.005696508567780256,
  0.005582868587225676,
  0.005560555029660463,
  0.005471393465995789,
  0.00543140945956111,
  0.005371361970901489,
  0.005345821846276522,
  0.005172608885914087,
  0.00497081084176898,
  0.004912454634904861,
  0.004871451295912266,
  0.004804256372153759,
  0.004735405091196299,
  0.004706101957708597,
  0.00468203192576766,
  0.004647127818316221,
  0.004544373601675034,
  0.004535902757197618
)

for ($i = 0; $i -lt $vars.Length; $i++) {
  $row = $i + 2
  $ws.Cells.Item($row, 2).Value = $vars[$i]
  $ws.Cells.Item($row, 3).Value = $scores[$i]
}